$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = "IT Project managment"
$find.Replacement.ClearFormatting()
$find.Replacement.Text = "IT Project management"
$find.Execute(
    $find.Text,      # FindText
    $true,           # MatchCase
    $true,           # MatchWholeWord
    $false,          # MatchWildcards
    $false,          # MatchSoundsLike
    $false,          # MatchAllWordForms
    $true,           # Forward
    1,               # Wrap (wdFindContinue)
    $false,          # Format
    $find.Replacement.Text, # ReplaceWith
    2                # Replace (wdReplaceAll)
) | Out-Null

Write-Output "done"
